# Issue with OTP: shorten the child-sub-category delete confirmation message,
# add new language keys (delete_child_subcat_label, phone/OTP verification
# strings) as new rows at the bottom of the language-variable sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: the long confirmation message is replaced with a shorter one.
$ws.Range("C5").Value = "Do you want to delete this child sub category?"
# The row no longer needs the manual 30pt height - let Excel re-fit it now
# that the text is short enough to wrap within the default row height.
$ws.Rows.Item(5).AutoFit()

# New row 27: delete_child_subcat_label / "Delete Child Sub Category?"
# B27 gets a distinct font/alignment (dark gray text, left aligned, wrap).
$ws.Range("A27").Value = "en"
$ws.Range("B27").Value = "delete_child_subcat_label"
$ws.Range("B27").Font.Color = 2696481
$ws.Range("B27").HorizontalAlignment = -4131
$ws.Range("B27").WrapText = $true
$ws.Range("C27").Value = "Delete Child Sub Category?"

# New row 28: err_user_phone_verification
$ws.Range("A28").Value = "en"
$ws.Range("B28").Value = "err_user_phone_verification"
$ws.Range("C28").Value = "Can not verify mobile phone number, please try again"

# New row 29: success_user_mobile_verify
$ws.Range("A29").Value = "en"
$ws.Range("B29").Value = "success_user_mobile_verify"
$ws.Range("C29").Value = "User mobile verified successfully"
